$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: the centred title at the very top of the letter currently spells
# out "CRIART CRIACOES PROMOCIONAIS EIRELE" across three runs ("...PROMOCIO"
# + "NAIS" + " EIRELE"). Revert it to read "...PROMOCIOANSI EIRELE" (the
# same garbled wording used later on in the body of the letter) by swapping
# the middle run's text from "NAIS" to "ANSI". The edit stays fully inside
# that run's own boundaries, so the (now textually adjacent, identically
# formatted) runs collapse into a single run, just like the original
# (pre-fix) document had it.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(1).Range
$null = $p1.Find.Execute("NAIS", $true, $false, $false, $false, $false, $true, 1, $false, "ANSI", 2)

# ---------------------------------------------------------------------------
# Edit 2: further down, in the body paragraph, "estabelecida na Rua Itambi"
# loses the "ua" of "Rua" (back to "estabelecida na R Itambi"). Delete just
# those two characters, which make up a whole run by themselves, so the
# surrounding runs collapse the same way the original document had them.
# ---------------------------------------------------------------------------
$p11 = $d.Paragraphs(11).Range
$idx = $p11.Text.IndexOf("Rua ")
$uaStart = $p11.Start + $idx + 1
$uaEnd = $uaStart + 2
$uaRange = $d.Range($uaStart, $uaEnd)
$uaRange.Delete()
